$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DRE")

# Row 17 - Net Profit Margin
$ws.Range("D17").Value = 0.2265
$ws.Range("E17").Value = 0.3971
$ws.Range("F17").Value = 0.4315
$ws.Range("G17").Value = 0.4405

# Row 18 - Free Cash Flow Margin
$ws.Range("B18").Value = -0.3898
$ws.Range("C18").Value = -0.4972
$ws.Range("D18").Value = -0.655
$ws.Range("E18").Value = 0.0908
$ws.Range("F18").Value = -0.2969
$ws.Range("G18").Value = -0.3476

# Row 19 - EBITDA (B19 was an empty inline string, now becomes a number)
$ws.Range("B19").Value = 836493000.0

# Row 21 - EPS (Diluted, from Cont. Ops)
$ws.Range("B21").Value = 0.96
$ws.Range("D21").Value = 0.5865
$ws.Range("E21").Value = 1.0165
$ws.Range("F21").Value = 1.1065

# Row 25 - EPS (Basic, from Continuous Ops)
$ws.Range("B25").Value = 0.97
$ws.Range("D25").Value = 0.5853
$ws.Range("E25").Value = 1.0153
$ws.Range("F25").Value = 1.1053

# Row 29 - EBITDA Margin
$ws.Range("B29").Value = 0.792

# Row 30 - Operating Cash Flow Margin
$ws.Range("B30").Value = 0.5709
$ws.Range("D30").Value = 0.5657
$ws.Range("E30").Value = 0.5437
$ws.Range("F30").Value = 0.5409
$ws.Range("G30").Value = 0.5195
